$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @(
    "0805053472082",
    "9405051381087",
    "9105050484185",
    "0805051084186",
    "9405054015187",
    "96050533240810",
    "9605052815088",
    "9605051240189",
    "9605050620183",
    "96050516500810",
    "86050540360810",
    "8605051812087",
    "86050505900810",
    "8605052718085",
    "8605051876082",
    "9605052337182",
    "9605054854085",
    "91050525460810",
    "9105050846185",
    "8805052747189",
    "8805050407083",
    "8805052090184",
    "8805050646086",
    "8805053064089",
    "8805054951185"
)

$startRow = 350
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $newValues[$i]
    $cell.NumberFormat = "@"
}
